$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the time_taken timestamps in column F (rows 2-86) of the "data" sheet
$timestamps = @(
    "2021-10-05 14:22:10.534829",
    "2021-10-05 14:22:10.534837",
    "2021-10-05 14:22:10.534840",
    "2021-10-05 14:22:10.534843",
    "2021-10-05 14:22:10.534846",
    "2021-10-05 14:22:10.534848",
    "2021-10-05 14:22:10.534851",
    "2021-10-05 14:22:10.534853",
    "2021-10-05 14:22:10.534856",
    "2021-10-05 14:22:10.534859",
    "2021-10-05 14:22:10.534862",
    "2021-10-05 14:22:10.534864",
    "2021-10-05 14:22:10.534866",
    "2021-10-05 14:22:10.534869",
    "2021-10-05 14:22:10.534872",
    "2021-10-05 14:22:10.534874",
    "2021-10-05 14:22:10.534877",
    "2021-10-05 14:22:10.534880",
    "2021-10-05 14:22:10.534882",
    "2021-10-05 14:22:10.534884",
    "2021-10-05 14:22:10.534887",
    "2021-10-05 14:22:10.534889",
    "2021-10-05 14:22:10.534892",
    "2021-10-05 14:22:10.534894",
    "2021-10-05 14:22:10.534897",
    "2021-10-05 14:22:10.534900",
    "2021-10-05 14:22:10.534902",
    "2021-10-05 14:22:10.534905",
    "2021-10-05 14:22:10.534907",
    "2021-10-05 14:22:10.534909",
    "2021-10-05 14:22:10.534912",
    "2021-10-05 14:22:10.534914",
    "2021-10-05 14:22:10.534917",
    "2021-10-05 14:22:10.534920",
    "2021-10-05 14:22:10.534922",
    "2021-10-05 14:22:10.534925",
    "2021-10-05 14:22:10.534927",
    "2021-10-05 14:22:10.534929",
    "2021-10-05 14:22:10.534932",
    "2021-10-05 14:22:10.534935",
    "2021-10-05 14:22:10.534938",
    "2021-10-05 14:22:10.534940",
    "2021-10-05 14:22:10.534943",
    "2021-10-05 14:22:10.534945",
    "2021-10-05 14:22:10.534947",
    "2021-10-05 14:22:10.534950",
    "2021-10-05 14:22:10.534952",
    "2021-10-05 14:22:10.534955",
    "2021-10-05 14:22:10.534957",
    "2021-10-05 14:22:10.534960",
    "2021-10-05 14:22:10.534962",
    "2021-10-05 14:22:10.534964",
    "2021-10-05 14:22:10.534967",
    "2021-10-05 14:22:10.534970",
    "2021-10-05 14:22:10.534972",
    "2021-10-05 14:22:10.534975",
    "2021-10-05 14:22:10.534977",
    "2021-10-05 14:22:10.534980",
    "2021-10-05 14:22:10.534982",
    "2021-10-05 14:22:10.534985",
    "2021-10-05 14:22:10.534987",
    "2021-10-05 14:22:10.534990",
    "2021-10-05 14:22:10.534992",
    "2021-10-05 14:22:10.534995",
    "2021-10-05 14:22:10.534998",
    "2021-10-05 14:22:10.535001",
    "2021-10-05 14:22:10.535004",
    "2021-10-05 14:22:10.535006",
    "2021-10-05 14:22:10.535009",
    "2021-10-05 14:22:10.535011",
    "2021-10-05 14:22:10.535014",
    "2021-10-05 14:22:10.535016",
    "2021-10-05 14:22:10.535019",
    "2021-10-05 14:22:10.535021",
    "2021-10-05 14:22:10.535024",
    "2021-10-05 14:22:10.535026",
    "2021-10-05 14:22:10.535031",
    "2021-10-05 14:22:10.535034",
    "2021-10-05 14:22:10.535037",
    "2021-10-05 14:22:10.535039",
    "2021-10-05 14:22:10.535041",
    "2021-10-05 14:22:10.535044",
    "2021-10-05 14:22:10.535046",
    "2021-10-05 14:22:10.535049",
    "2021-10-05 14:22:10.535052"
)
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

# Add a new "metadata" worksheet positioned right after "data"
$metaSheet = $wb.Worksheets.Add($null, $ws)
$metaSheet.Name = "metadata"

# Copy the header style (bold, centered, bordered) used by row 1 / column A
# of the "data" sheet so the new sheet's styled cells reuse the same style index.
$ws.Range("B1").Copy($metaSheet.Range("B1"))
$ws.Range("B1").Copy($metaSheet.Range("C1"))
$ws.Range("B1").Copy($metaSheet.Range("D1"))
$ws.Range("B1").Copy($metaSheet.Range("E1"))
$ws.Range("B1").Copy($metaSheet.Range("F1"))
$ws.Range("B1").Copy($metaSheet.Range("G1"))
$ws.Range("A2").Copy($metaSheet.Range("A2"))

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Paroxysmal central nervous system disorders"
$metaSheet.Range("C2").Value = 541

# "1.17" must stay a text string (not coerce to the number 1.17) while keeping
# the cell on the default (unstyled) format, so temporarily mark it as Text,
# assign the value, then restore the Normal cell style.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.17"
$metaSheet.Range("D2").Style = "Normal"

$metaSheet.Range("E2").Value = "2021-08-04T12:56:37.015099Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:22:10.531633"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/541/?format=json"
